$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Förändrad" (column C) changed from 2023-09-23 (45192) to 2023-10-03 (45202)
#    for every existing data row (2..323).
$ws.Range("C2:C323").Value = 45202

# 2) Row 323 gains an explicit row height (ht="15" customHeight="1"),
#    matching the rest of the data rows.
$ws.Rows.Item(323).RowHeight = 15

# 3) Two new rows are appended: 324 and 325.

# Row 324
$ws.Cells.Item(324, 1).Value = "A 46784-2023"
$ws.Cells.Item(324, 2).Value = 45198
$ws.Cells.Item(324, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(324, 3).Value = 45202
$ws.Cells.Item(324, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(324, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(324, 5).Value = "ANEBY"
$ws.Cells.Item(324, 6).Value = "Sveaskog"
$ws.Cells.Item(324, 7).Value = 0.9
$ws.Range("H324:Q324").Value = 0
$ws.Range("R324").WrapText = $true
$ws.Rows.Item(324).RowHeight = 15

# Row 325
$ws.Cells.Item(325, 1).Value = "A 46783-2023"
$ws.Cells.Item(325, 2).Value = 45198
$ws.Cells.Item(325, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(325, 3).Value = 45202
$ws.Cells.Item(325, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(325, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(325, 5).Value = "ANEBY"
$ws.Cells.Item(325, 6).Value = "Sveaskog"
$ws.Cells.Item(325, 7).Value = 1.6
$ws.Range("H325:Q325").Value = 0
$ws.Range("R325").WrapText = $true
